# ------------------------------------------------------------------
# Adds a new "ELC_IND_FIN" (Electricity - Industry Final Energy)
# commodity and its related "ELC_IND_FIN_DEM" demand process to the
# VEDA/TIMES model workbook, mirroring the existing ELC_FIN /
# ELC_FIN_DEM rows across the SEC_Comm, SEC_Processes,
# FINAL_DEMAD_PRC and DEMAND sheets.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsComm    = $wb.Worksheets.Item("SEC_Comm")
$wsProc    = $wb.Worksheets.Item("SEC_Processes")
$wsFinal   = $wb.Worksheets.Item("FINAL_DEMAD_PRC")
$wsDemand  = $wb.Worksheets.Item("DEMAND")

# --------------------------------------------------------------
# 1. SEC_Comm: define the new "ELC_IND_FIN" commodity on row 10
#    (B10:G10 mirrors the existing ELC_FIN definition on row 9)
# --------------------------------------------------------------
$wsComm.Range("B10").Value2 = "DEM"
$wsComm.Range("C10").Value2 = "ELC_IND_FIN"
$wsComm.Range("D10").Value2 = "Electricity - Industry Final Energy"
$wsComm.Range("E10").Value2 = "PJ"
$wsComm.Range("G10").Value2 = "DAYNITE"

# G10 picks up the same highlighted style already used by G9
$wsComm.Range("G9").Copy() | Out-Null
$wsComm.Range("G10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --------------------------------------------------------------
# 2. SEC_Processes: rename the existing demand process label and
#    add the new "ELC_IND_FIN_DEM" process on row 9
# --------------------------------------------------------------
$wsProc.Range("B9").Value2 = "DMD"
$wsProc.Range("C9").Value2 = "NL"
$wsProc.Range("D9").Value2 = "ELC_IND_FIN_DEM"

$wsProc.Range("E8").Value2 = "Electricity - Industry Final Energy Demand"

$wsProc.Range("E9").Value2 = "Electricity - Final Energy Demand"
$wsProc.Range("F9").Value2 = "PJ"
$wsProc.Range("G9").Value2 = "PJa"
$wsProc.Range("H9").Value2 = "DAYNITE"
$wsProc.Range("J9").Value2 = "NO"

# Match the formatting already used on the row above for the cells
# that carry the shaded / highlighted style
$wsProc.Range("D8").Copy() | Out-Null
$wsProc.Range("E9").PasteSpecial(-4122) | Out-Null

$wsProc.Range("F8").Copy() | Out-Null
$wsProc.Range("F9").PasteSpecial(-4122) | Out-Null

$wsProc.Range("G8").Copy() | Out-Null
$wsProc.Range("G9").PasteSpecial(-4122) | Out-Null

$wsProc.Range("H8").Copy() | Out-Null
$wsProc.Range("H9").PasteSpecial(-4122) | Out-Null

$wsProc.Range("I8").Copy() | Out-Null
$wsProc.Range("I9").PasteSpecial(-4122) | Out-Null

$wsProc.Range("J8").Copy() | Out-Null
$wsProc.Range("J9").PasteSpecial(-4122) | Out-Null

# --------------------------------------------------------------
# 3. FINAL_DEMAD_PRC: add the ELC_IND_FIN_DEM process block
#    (rows 12-16), mirroring the ELC_FIN_DEM block on rows 8-11
# --------------------------------------------------------------
$wsFinal.Range("B12").Formula = "=SEC_Processes!D9"
$wsFinal.Range("D12").Formula = "=SEC_Comm!C11"
$wsFinal.Range("F12").Value2 = 1
$wsFinal.Range("G12").Value2 = 1

$wsFinal.Range("D13").Formula = "=SEC_Comm!C8"
$wsFinal.Range("D14").Value2 = "ELC_GRID_RES"
$wsFinal.Range("D15").Value2 = "ELC_IND_RES"
$wsFinal.Range("E16").Formula = "=SEC_Comm!C10"

# Apply the shaded style (already used by D9:D11) to the new D column
# cells, and extend it down onto the previously-blank D17 cell too
$wsFinal.Range("D9").Copy() | Out-Null
$wsFinal.Range("D12:D17").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# F8/G8 style is reused for both F12 and G12
$wsFinal.Range("F8").Copy() | Out-Null
$wsFinal.Range("F12").PasteSpecial(-4122) | Out-Null
$wsFinal.Range("G12").PasteSpecial(-4122) | Out-Null

# --------------------------------------------------------------
# 4. DEMAND: fill in the ELC_IND_FIN annual demand values (row 9)
# --------------------------------------------------------------
$wsDemand.Range("B9").Formula = "=SEC_Comm!C10"
$wsDemand.Range("C9").Value2 = 100
$wsDemand.Range("D9").Value2 = 12
$wsDemand.Range("E9").Value2 = 150
$wsDemand.Range("F9").Value2 = 200
$wsDemand.Range("G9").Value2 = 300
$wsDemand.Range("H9").Value2 = 400
$wsDemand.Range("I9").Value2 = 500

# --------------------------------------------------------------
# 5. Recalculate so every cached formula value reflects the new data
# --------------------------------------------------------------
$excel.CalculateFullRebuild()

# --------------------------------------------------------------
# 6. Leave the view focused on FINAL_DEMAD_PRC, matching the final
#    state captured in the saved workbook
# --------------------------------------------------------------
$wsFinal.Activate()
$wsFinal.Range("G19").Select() | Out-Null
